$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record as row 23, pushing the existing rows 23-24 down to 24-25.
$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44516
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112026
$ws.Range("G23").Value = "Haba"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 7500
$ws.Range("N23").Value = "`$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 300
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"

# Update values that changed on what is now row 24 (formerly row 23).
$ws.Range("J24").Value = 100
